# Documentation.xlsx - add a new catalog row for the "Data Model Links" asset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns are: A=ID, B=File Name, C=Description, D=Directory
# Write File Name / Description / Directory / ID in this order so the new
# shared-string table entries land in the same order as the source edit.
$ws.Range("B2").Value = "Data Model Links.png"
$ws.Range("C2").Value = "Data Model Links"
$ws.Range("D2").Value = "data/Multimedia_Data/Documentation/"
$ws.Range("A2").Value = "DOC_001"

# Widen the data columns to fit the new content.
$ws.Columns.Item(2).ColumnWidth = 18.33203125
$ws.Columns.Item(3).ColumnWidth = 25.6640625
$ws.Columns.Item(4).ColumnWidth = 34.1640625

$ws.Range("A3").Select()
